# Updates the "cryptos" price/volume snapshot (GitHub Actions scheduled refresh).
# Columns: B=Coin, C=Link, D=Price, E=Volume(1h). Rows 46/47 also swap the
# Bittensor/Maker entries (their rank order changed in this refresh).
#
# Price strings in column D must stay TEXT (they use dot-grouped thousands,
# e.g. "66.977.85", and some need exact trailing-zero formatting like
# "36.90"), so plain-numeric-looking values are written with a leading
# apostrophe (quote-prefix) to stop Excel's COM layer from auto-coercing
# them to doubles, then the cell style is reset to "Normal" so the
# quote-prefix formatting doesn't linger on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.977.85'
$ws.Range("E2").Value = '  +2.48%  '
$ws.Range("D3").Value = '3.103.84'
$ws.Range("E3").Value = '  +5.35%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '''579.55'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.62%  '
$ws.Range("D6").Value = '''171.82'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.55%  '
$ws.Range("D7").Value = '''0.999'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").Value = '3.095.61'
$ws.Range("E8").Value = '  +5.20%  '
$ws.Range("E9").Value = '  +1.42%  '
$ws.Range("D10").Value = '''6.59'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -1.71%  '
$ws.Range("D11").Value = '''0.156'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +3.49%  '
$ws.Range("D12").Value = '''0.483'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +5.40%  '
$ws.Range("E13").Value = '  +2.18%  '
$ws.Range("D14").Value = '''36.90'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.94%  '
$ws.Range("E15").Value = '  -0.48%  '
$ws.Range("D16").Value = '3.610.41'
$ws.Range("E16").Value = '  +5.13%  '
$ws.Range("D17").Value = '66.914.34'
$ws.Range("E17").Value = '  +2.36%  '
$ws.Range("D18").Value = '''7.21'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.79%  '
$ws.Range("D19").Value = '3.094.67'
$ws.Range("E19").Value = '  +5.01%  '
$ws.Range("D20").Value = '''16.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.85%  '
$ws.Range("D21").Value = '''480.73'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.27%  '
$ws.Range("D22").Value = '''0.717'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.93%  '
$ws.Range("D23").Value = '''7.52'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +3.19%  '
$ws.Range("D24").Value = '''83.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.66%  '
$ws.Range("D25").Value = '''13.06'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +7.46%  '
$ws.Range("D26").Value = '''2.34'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.50%  '
$ws.Range("D27").Value = '''10.08'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.90%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("D29").Value = '''8.01'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.17%  '
$ws.Range("E30").Value = '  -2.79%  '
$ws.Range("D31").Value = '''2.68'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +3.76%  '
$ws.Range("D32").Value = '''0.0000102'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +2.09%  '
$ws.Range("D33").Value = '''28.43'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +4.27%  '
$ws.Range("E34").Value = '  +1.92%  '
$ws.Range("E35").Value = '  +0.01%  '
$ws.Range("D36").Value = '''1.01'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.59%  '
$ws.Range("D37").Value = '''5.91'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +3.25%  '
$ws.Range("D38").Value = '''47.86'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +5.66%  '
$ws.Range("D39").Value = '''2.12'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +10.09%  '
$ws.Range("D40").Value = '''50.19'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.55%  '
$ws.Range("E41").Value = '  +4.77%  '
$ws.Range("E42").Value = '  +0.62%  '
$ws.Range("D43").Value = '''8.67'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.49%  '
$ws.Range("D44").Value = '''2.79'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.69%  '
$ws.Range("D45").Value = '''0.0360'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +2.67%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '2.816.40'
$ws.Range("E46").Value = '  +5.03%  '
$ws.Range("B47").Value = 'Bittensor'
$ws.Range("C47").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D47").Value = '''383.06'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.35%  '
$ws.Range("D48").Value = '''134.63'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.15%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("D50").Value = '''24.86'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.45%  '
$ws.Range("E51").Value = '  +2.72%  '
